$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last column (F) entirely - data now only spans A:E
$ws.Range("F1:F9").Delete()

# Update header row (B1:E1)
$ws.Range("B1").Value = "IBOVESPA"
$ws.Range("C1").Value = "IMA-B"
$ws.Range("D1").Value = "SCBTG"
$ws.Range("E1").Value = "S&P 500"

# Row 2 - Mean
$ws.Range("B2").Value = 114014.0912738215
$ws.Range("C2").Value = 8858.345973007021
$ws.Range("D2").Value = 7.977432296890672
$ws.Range("E2").Value = 246.6238513971916

# Row 3 - Median
$ws.Range("B3").Value = 112764
$ws.Range("C3").Value = 8569.517189
$ws.Range("D3").Value = 8.039999999999999
$ws.Range("E3").Value = 238.789993

# Row 4 - Standard Deviation
$ws.Range("B4").Value = 8847.243905469244
$ws.Range("C4").Value = 700.844902350853
$ws.Range("D4").Value = 0.745963949075289
$ws.Range("E4").Value = 25.03584273367861

# Row 5 - Kurtosis
$ws.Range("B5").Value = -0.711850269504962
$ws.Range("C5").Value = -1.451583380115703
$ws.Range("D5").Value = 0.9528817755822936
$ws.Range("E5").Value = -0.5655566533230973

# Row 6 - Skewness
$ws.Range("B6").Value = 0.3556782972634858
$ws.Range("C6").Value = 0.210039751158321
$ws.Range("D6").Value = 0.5232456205924026
$ws.Range("E6").Value = 0.7397614525688554

# Row 7 - Fishers Information
$ws.Range("B7").Value = 0.1003180021900583
$ws.Range("C7").Value = 0.1089998067752257
$ws.Range("D7").Value = 0.1127085498772781
$ws.Range("E7").Value = 0.09133581486054189

# Row 8 - MIEE
$ws.Range("B8").Value = 0.2441202301271639
$ws.Range("C8").Value = 0.2861911670287034
$ws.Range("D8").Value = 0.2740876642375576
$ws.Range("E8").Value = 0.2301067906714883

# Row 9 - Permutation Entropy
$ws.Range("B9").Value = 0.8561977720628944
$ws.Range("C9").Value = 0.8228086397465223
$ws.Range("D9").Value = 0.8386208856397205
$ws.Range("E9").Value = 0.8612290241890536
